# Update TIRAFONDOS price list: refresh the date stamp in A1 and the
# unit-price values in column D for the 3/16", 1/4" and 5/16" tirafondos
# sections (rows 32-40, 46-61, 67-81 and 87-96) on the first worksheet
# ("Hoja1").

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$updates = @{
    "A1" = 45436

    "D32" = 2236.3
    "D33" = 2673
    "D34" = 2992
    "D35" = 3223
    "D36" = 3410
    "D37" = 3829.1
    "D38" = 4213
    "D39" = 4510
    "D40" = 4752

    "D46" = 2728
    "D47" = 3047
    "D48" = 3344
    "D49" = 3894
    "D50" = 4236.1
    "D51" = 4719
    "D52" = 5297.6
    "D53" = 5522
    "D54" = 6094
    "D55" = 6809
    "D56" = 7393.1
    "D57" = 8305
    "D58" = 9086
    "D59" = 9955
    "D60" = 11348.823
    "D61" = 12044.74

    "D67" = 5049
    "D68" = 5758.5
    "D69" = 6204
    "D70" = 6303
    "D71" = 7117
    "D72" = 7667
    "D73" = 8538.200000000001
    "D74" = 9168.5
    "D75" = 10642.5
    "D76" = 11572
    "D77" = 12881
    "D78" = 14179
    "D79" = 15730
    "D80" = 21703
    "D81" = 24178

    "D87" = 8954
    "D88" = 9240
    "D89" = 10824
    "D90" = 11858
    "D91" = 13090
    "D92" = 15994
    "D93" = 16643
    "D94" = 18480
    "D95" = 20240
    "D96" = 22385
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
